$wb = $excel.ActiveWorkbook

# --- Sheet "Metrics": move selection only ---
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsMetrics.Range("C15").Select()

# --- Sheet "csdjzqs": move selection, clear out row 10's data (keep B10/C10 styled-empty) ---
$wsCsdjzqs = $wb.Worksheets.Item("csdjzqs")
$wsCsdjzqs.Range("A10").ClearContents()
$wsCsdjzqs.Range("B10").ClearContents()
$wsCsdjzqs.Range("C10").ClearContents()
$wsCsdjzqs.Range("H42").Select()

# --- Sheet "ndzsrqs": move selection, delete entire row 10 ---
$wsNdzsrqs = $wb.Worksheets.Item("ndzsrqs")
$wsNdzsrqs.Rows(10).Delete()
$wsNdzsrqs.Range("G8").Select()

# --- Sheet "bksr": update values, move selection/top row, keep as active sheet ---
$wsBksr = $wb.Worksheets.Item("bksr")
$wsBksr.Range("C5").Value = 68818
$wsBksr.Range("C15").Value = 126047.14
$wsBksr.Range("C25").Value = 25675
$wsBksr.Activate()
$wsBksr.Range("A16").Select()
$wsBksr.Range("F24:G24").Select()
